$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 0
$ws1.Range("F8").Value = 147
$ws1.Range("F9").Value = 65
$ws1.Range("F10").Value = 512

# Sheet "全部类型" (fourth sheet)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 119
$ws4.Range("F9").Value = 65
